# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets to reflect that the
# two source files have been handed back and are in sync with en-US.

$wb = $excel.ActiveWorkbook

# Target (source) file -> handoff file url, per language, used to populate the
# "Latest Target File" / "Latest Handback File" columns and their hyperlinks.
$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/29279770b0f141116d0a4a8b35afad74239d05a6/e2e/0f776044-d944-43a8-a3ba-88fe1f128e4a.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/29279770b0f141116d0a4a8b35afad74239d05a6/e2e/2d156d40-711f-4eff-b23d-464f0877684f.md"

$zhXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494236234cd56fb3d31584e5c23ac78daab16b22/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0f776044-d944-43a8-a3ba-88fe1f128e4a.f56e2ac6fdde5b656d6547081c795cddca49f986.zh-cn.xlf"
$zhXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494236234cd56fb3d31584e5c23ac78daab16b22/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2d156d40-711f-4eff-b23d-464f0877684f.931b21affdae1bebb75d27ef4f375b649b33c0cf.zh-cn.xlf"

$deXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40974ceb7f064739bccf64e465dda15a2d6918e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0f776044-d944-43a8-a3ba-88fe1f128e4a.f56e2ac6fdde5b656d6547081c795cddca49f986.de-de.xlf"
$deXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40974ceb7f064739bccf64e465dda15a2d6918e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2d156d40-711f-4eff-b23d-464f0877684f.931b21affdae1bebb75d27ef4f375b649b33c0cf.de-de.xlf"

$mdName1 = "0f776044-d944-43a8-a3ba-88fe1f128e4a.md"
$mdName2 = "2d156d40-711f-4eff-b23d-464f0877684f.md"
$zhXlfName1 = "0f776044-d944-43a8-a3ba-88fe1f128e4a.f56e2ac6fdde5b656d6547081c795cddca49f986.zh-cn.xlf"
$zhXlfName2 = "2d156d40-711f-4eff-b23d-464f0877684f.931b21affdae1bebb75d27ef4f375b649b33c0cf.zh-cn.xlf"
$deXlfName1 = "0f776044-d944-43a8-a3ba-88fe1f128e4a.f56e2ac6fdde5b656d6547081c795cddca49f986.de-de.xlf"
$deXlfName2 = "2d156d40-711f-4eff-b23d-464f0877684f.931b21affdae1bebb75d27ef4f375b649b33c0cf.de-de.xlf"

$statusText = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (B) - mark both rows as handed back
$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

# New "Latest Target File" (E) / "Latest Handback File" (F) hyperlinked cells
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl1, "", "", $mdName1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl1, "", "", $zhXlfName1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl2, "", "", $mdName2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl2, "", "", $zhXlfName2) | Out-Null

# Latest Handback DateTime (G) - now filled in with actual handback timestamps
$wsZh.Range("G2").Value = "2016-02-15 08:59:18"
$wsZh.Range("G3").Value = "2016-02-15 08:59:18"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl1, "", "", $mdName1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl1, "", "", $deXlfName1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl2, "", "", $mdName2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl2, "", "", $deXlfName2) | Out-Null

$wsDe.Range("G2").Value = "2016-02-15 08:59:45"
$wsDe.Range("G3").Value = "2016-02-15 08:59:45"
